$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.826.58"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.751.80"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.91"
$ws.Range("E5").Value = "  -1.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.19"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("E8").Value = "  -1.85%  "
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.381"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.61"
$ws.Range("E12").Value = "  -17.11%  "
$ws.Range("D13").Value = "3.236.96"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.41"
$ws.Range("E14").Value = "  -2.04%  "
$ws.Range("D15").Value = "63.460.20"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000149"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "2.754.69"
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.10"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  -2.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.92"
$ws.Range("E20").Value = "  -2.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("E21").Value = "  -4.59%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.535"
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.14"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("E25").Value = "  -1.62%  "
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.38"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "0.0₃0899"
$ws.Range("E28").Value = "  -0.79%  "
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.19"
$ws.Range("E31").Value = "  -3.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.10"
$ws.Range("E33").Value = "  -2.75%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.84"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.977"
$ws.Range("E38").Value = "  -4.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.15"
$ws.Range("E39").Value = "  +4.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.12"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "326.48"
$ws.Range("E41").Value = "  -4.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.89"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.33"
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0583"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.22"
$ws.Range("E45").Value = "  -4.14%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.93"
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0252"
$ws.Range("E47").Value = "  -2.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.622"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.100"
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("E50").Value = "  +0.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.44%  "
